# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Tue May 14 16:56:32 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (column D holds numeric-looking
# strings such as "61.327.87" or "31.90" that must stay text, not be
# coerced into numbers / lose trailing zeros). Force text format, assign,
# then restore the cell to the "Normal" style so no stray formatting is
# left behind.
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "61.327.87"
$ws.Range("E2").Value = "  -2.92%  "
Set-TextValue "D3" "2.889.63"
$ws.Range("E3").Value = "  -2.77%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "566.01"
$ws.Range("E5").Value = "  -5.02%  "
Set-TextValue "D6" "143.11"
$ws.Range("E6").Value = "  -4.24%  "
$ws.Range("E7").Value = "  +0.21%  "
Set-TextValue "D8" "0.506"
$ws.Range("E8").Value = "  -0.78%  "
Set-TextValue "D9" "2.888.16"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("E10").Value = "  -8.72%  "
Set-TextValue "D11" "0.146"
$ws.Range("E11").Value = "  -5.48%  "
$ws.Range("E12").Value = "  -3.19%  "
Set-TextValue "D13" "0.0000233"
$ws.Range("E13").Value = "  -3.73%  "
Set-TextValue "D14" "31.90"
$ws.Range("E14").Value = "  -4.13%  "
$ws.Range("E15").Value = "  -0.67%  "
Set-TextValue "D16" "3.367.60"
$ws.Range("E16").Value = "  -2.71%  "
Set-TextValue "D17" "61.337.10"
$ws.Range("E17").Value = "  -2.75%  "
Set-TextValue "D18" "6.62"
$ws.Range("E18").Value = "  -2.64%  "
Set-TextValue "D19" "2.886.54"
$ws.Range("E19").Value = "  -2.58%  "
Set-TextValue "D20" "432.16"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("E22").Value = "  -2.94%  "
Set-TextValue "D23" "6.85"
$ws.Range("E23").Value = "  -4.05%  "
Set-TextValue "D24" "79.24"
$ws.Range("E24").Value = "  -3.24%  "
Set-TextValue "D25" "11.80"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D26" "10.01"
$ws.Range("E26").Value = "  -12.52%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D27" "1.00"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -8.95%  "
$ws.Range("E29").Value = "  -4.22%  "
Set-TextValue "D30" "7.01"
$ws.Range("E30").Value = "  -4.36%  "
$ws.Range("E31").Value = "  -5.51%  "
$ws.Range("E32").Value = "  -8.92%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -3.44%  "
Set-TextValue "D35" "25.59"
$ws.Range("E35").Value = "  -4.70%  "
Set-TextValue "D36" "0.955"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("E37").Value = "  -4.55%  "
Set-TextValue "D38" "48.70"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("E39").Value = "  -13.85%  "
$ws.Range("E40").Value = "  -6.73%  "
Set-TextValue "D41" "8.25"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("E42").Value = "  -4.17%  "
Set-TextValue "D43" "39.69"
$ws.Range("E43").Value = "  -4.31%  "
Set-TextValue "D44" "0.268"
$ws.Range("E44").Value = "  -7.27%  "
Set-TextValue "D45" "2.686.72"
$ws.Range("E45").Value = "  -1.12%  "
Set-TextValue "D46" "133.99"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("E48").Value = "  -0.01%  "
Set-TextValue "D49" "339.38"
$ws.Range("E49").Value = "  -8.83%  "
$ws.Range("E50").Value = "  -2.75%  "
Set-TextValue "D51" "21.47"
$ws.Range("E51").Value = "  -7.70%  "
